$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '60.940.74'
$ws.Cells.Item(2, 5).Value = '  +0.28%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '2.350.65'
$ws.Cells.Item(3, 5).Value = '  -0.72%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.03%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '544.41'
$ws.Cells.Item(5, 5).Value = '  +0.30%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '136.82'
$ws.Cells.Item(6, 5).Value = '  -2.42%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  -0.03%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.525'
$ws.Cells.Item(8, 5).Value = '  -6.17%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '2.347.85'
$ws.Cells.Item(9, 5).Value = '  -0.74%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '0.105'
$ws.Cells.Item(10, 5).Value = '  +0.06%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  +1.75%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '5.30'
$ws.Cells.Item(12, 5).Value = '  -0.60%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '0.342'
$ws.Cells.Item(13, 5).Value = '  +0.14%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '24.67'
$ws.Cells.Item(14, 5).Value = '  -2.70%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '2.772.62'
$ws.Cells.Item(15, 5).Value = '  -0.75%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '60.798.35'
$ws.Cells.Item(16, 5).Value = '  +0.16%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '0.0000160'
$ws.Cells.Item(17, 5).Value = '  -1.96%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '2.343.09'
$ws.Cells.Item(18, 5).Value = '  -1.09%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '10.64'
$ws.Cells.Item(19, 5).Value = '  +0.98%  '

# Row 20
$ws.Cells.Item(20, 2).Value = 'Polkadot'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '4.12'
$ws.Cells.Item(20, 5).Value = '  +0.61%  '

# Row 21
$ws.Cells.Item(21, 2).Value = 'BitcoinCash'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '318.73'
$ws.Cells.Item(21, 5).Value = '  +0.67%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '6.54'
$ws.Cells.Item(22, 5).Value = '  -1.99%  '

# Row 23
$ws.Cells.Item(23, 5).Value = '  +0.11%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '63.36'
$ws.Cells.Item(24, 5).Value = '  +0.58%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '1.68'
$ws.Cells.Item(25, 5).Value = '  -6.36%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '8.38'
$ws.Cells.Item(26, 5).Value = '  +8.40%  '

# Row 27
$ws.Cells.Item(27, 5).Value = '  -0.03%  '

# Row 28
$ws.Cells.Item(28, 4).Value = '2.465.44'
$ws.Cells.Item(28, 5).Value = '  -0.74%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '7.94'
$ws.Cells.Item(29, 5).Value = '  -0.12%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '499.50'
$ws.Cells.Item(30, 5).Value = '  -3.25%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '1.37'
$ws.Cells.Item(31, 5).Value = '  -3.57%  '

# Row 32
$ws.Cells.Item(32, 2).Value = 'PEPE'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(32, 4).Value = '0.0₃0860'
$ws.Cells.Item(32, 5).Value = '  -7.28%  '

# Row 33
$ws.Cells.Item(33, 2).Value = 'Kaspa'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '0.145'
$ws.Cells.Item(33, 5).Value = '  +0.47%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  -2.09%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '1.49'
$ws.Cells.Item(35, 5).Value = '  -3.47%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '0.999'
$ws.Cells.Item(36, 5).Value = '  +0.03%  '

# Row 37
$ws.Cells.Item(37, 2).Value = 'NEARProtocol'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '4.59'
$ws.Cells.Item(37, 5).Value = '  -0.67%  '

# Row 38
$ws.Cells.Item(38, 2).Value = 'PolygonEcosystemToken'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '0.375'
$ws.Cells.Item(38, 5).Value = '  +0.54%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '18.44'
$ws.Cells.Item(39, 5).Value = '  +2.25%  '

# Row 40
$ws.Cells.Item(40, 2).Value = 'Stacks'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '1.83'
$ws.Cells.Item(40, 5).Value = '  +6.28%  '

# Row 41
$ws.Cells.Item(41, 2).Value = 'RenderToken'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '5.24'
$ws.Cells.Item(41, 5).Value = '  -3.59%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '143.14'
$ws.Cells.Item(42, 5).Value = '  +4.58%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '40.59'
$ws.Cells.Item(44, 5).Value = '  +1.22%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '142.95'
$ws.Cells.Item(45, 5).Value = '  +2.79%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '3.55'
$ws.Cells.Item(46, 5).Value = '  +0.73%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '2.04'
$ws.Cells.Item(47, 5).Value = '  -8.81%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '0.0518'
$ws.Cells.Item(48, 5).Value = '  +1.13%  '

# Row 49
$ws.Cells.Item(49, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '19.06'
$ws.Cells.Item(49, 5).Value = '  -6.71%  '

# Row 50
$ws.Cells.Item(50, 2).Value = 'Mantle'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '0.567'
$ws.Cells.Item(50, 5).Value = '  -1.23%  '

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '0.0899'
$ws.Cells.Item(51, 5).Value = '  -2.02%  '

Write-Output "Updated cryptos list"